# 10 icdc scripts for jenkins
# Update the FilesTab Cypher query (B4 on the "startup" sheet): drop the
# `File Type` and `Breed` columns from the RETURN clause.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = "MATCH (f:file)-->(parent)`r`n" + `
    "WITH DISTINCT f, parent`r`n" + `
    "MATCH (f)-[*]->(c:case)<--(demo:demographic)`r`n" + `
    " MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`r`n" + `
    "WHERE s.clinical_study_designation IN ['NCATS-COP01'] `r`n" + `
    "WITH DISTINCT f, parent, c, demo, diag, s`r`n" + `
    "RETURN coalesce(f.file_name, '') AS ``File Name``, `r`n" + `
    "        coalesce(labels(parent)[0], '') AS ``Association``,`r`n" + `
    "        coalesce(f.file_description, '') AS ``Description``,`r`n" + `
    "        coalesce(f.file_format, '') AS ``Format``,`r`n" + `
    "        coalesce(f.file_size, '') AS ``Size``,`r`n" + `
    "        coalesce(c.case_id, '') AS ``Case ID``, `r`n" + `
    "        coalesce(diag.disease_term,'') AS Diagnosis , `r`n" + `
    "        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newQuery

# Row 4 wraps this text at the default 14.5pt line height; with two fewer
# lines the autofit height shrinks from 232 to 203.
$ws.Rows.Item(4).RowHeight = 203

# Move the active selection / top-left scroll position down to row 4 (the
# cell the automation script was last working with).
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
